$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data corrections (pseudonymisation clean-up) ---
# Row 3 (Berta Brunner): civil status had been mistagged "L" -> should be "H"
$ws.Range("F3").Value = "Civil stat#H"
# Row 4 (house number on Clausiensteig): stray trailing "!" on the 4-char code -> "3c"
$ws.Range("K4").Value = "3c"

# --- Drop the EGID / EWID identifier columns entirely (PII-ish building/dwelling IDs) ---
# Everything to the right (VERMÖGEN, STEUERBARESEINKOMMEN, HASEL, HASSH, AMOUNT) shifts left.
$ws.Columns("N:O").Delete()

# --- Drop the three leaked/duplicate rows at the bottom (old rows 9-11) ---
$ws.Rows("9:11").Delete()

# --- Selection / view bookkeeping to match the saved state ---
$ws.Range("R9").Select()

# --- Window geometry as last recorded by Excel on save (best effort) ---
$excel.ActiveWindow.Left = 29400
$excel.ActiveWindow.Top = -7260
$excel.ActiveWindow.Width = 34600
$excel.ActiveWindow.Height = 21100
